$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.474.90"
$ws.Range("E2").Value = "  +0.40%  "
$ws.Range("D3").Value = "1.574.92"
$ws.Range("E3").Value = "  +0.91%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "288.22"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.80%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3705"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.82%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "47.73"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.92%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3329"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.151"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.37%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07568"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.49%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.04%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.85"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.960"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.86%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.950"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.26%  "
$ws.Range("D16").Value = "1.568.65"
$ws.Range("E16").Value = "  +0.50%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001124"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.32%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "88.31"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.40%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06733"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.27%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.409"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.59%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9998"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.55"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.49%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.04"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.15%  "
$ws.Range("D24").Value = "22.467.14"
$ws.Range("E24").Value = "  +0.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.385"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.638"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "151.48"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.63%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.68"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.71%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.993"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.37%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.52"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.19%  "
$ws.Range("D31").Value = "1.745.15"
$ws.Range("E31").Value = "  +0.60%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.095"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.96%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.118"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.45%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.986"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.22%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.867"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.26%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08351"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.45%  "
$ws.Range("E37").Value = "  +4.28%  "
$ws.Range("E38").Value = "  +1.80%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06402"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.05%  "
$ws.Range("E40").Value = "  +0.38%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.371"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.32%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.50"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.70%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6291"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.15%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.14"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.58%  "
$ws.Range("B45").Value = "Frax"
$ws.Range("C45").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.000"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6121"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +7.15%  "
$ws.Range("E47").Value = "  +0.64%  "
$ws.Range("E48").Value = "  +3.07%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "125.49"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.08%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.211"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.27%  "
$ws.Range("E51").Value = "  +0.12%  "
